$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.153.02"
$ws.Range("D3").Value = "1.656.82"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("D5").Value = "'218.85"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").Value = "'0.5244"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("E8").Value = "  +1.58%  "
$ws.Range("D9").Value = "'0.06357"
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("D10").Value = "'20.58"
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("D11").Value = "'0.07690"
$ws.Range("E11").Value = "  -1.52%  "
$ws.Range("D12").Value = "'4.616"
$ws.Range("E12").Value = "  +2.77%  "
$ws.Range("D13").Value = "1.659.55"
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("D14").Value = "1.885.09"
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("D15").Value = "'0.5625"
$ws.Range("E15").Value = "  +1.38%  "
$ws.Range("D16").Value = "0.0₅8200"
$ws.Range("E16").Value = "  +2.30%  "
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").Value = "26.151.66"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").Value = "'4.664"
$ws.Range("E20").Value = "  +0.59%  "
$ws.Range("D21").Value = "'10.57"
$ws.Range("E21").Value = "  +4.62%  "
$ws.Range("D22").Value = "'193.41"
$ws.Range("E22").Value = "  -1.21%  "
$ws.Range("D23").Value = "'5.960"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").Value = "'7.272"
$ws.Range("E27").Value = "  +1.71%  "
$ws.Range("D28").Value = "'15.98"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").Value = "'1.518"
$ws.Range("E29").Value = "  +2.10%  "
$ws.Range("D30").Value = "'0.05488"
$ws.Range("E30").Value = "  -3.89%  "
$ws.Range("D31").Value = "'1.273"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("D32").Value = "'3.468"
$ws.Range("E32").Value = "  -0.45%  "
$ws.Range("D33").Value = "'3.369"
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("D34").Value = "'1.564"
$ws.Range("E34").Value = "  -1.28%  "
$ws.Range("D35").Value = "'0.9533"
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("E36").Value = "  -0.86%  "
$ws.Range("D37").Value = "'2.403"
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("D38").Value = "'0.5697"
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("D39").Value = "'0.01590"
$ws.Range("E39").Value = "  -0.36%  "
$ws.Range("D40").Value = "'5.881"
$ws.Range("E40").Value = "  -1.12%  "
$ws.Range("D42").Value = "1.026.58"
$ws.Range("E42").Value = "  -3.64%  "
$ws.Range("D43").Value = "'0.8297"
$ws.Range("E43").Value = "  -1.89%  "
$ws.Range("D44").Value = "'101.31"
$ws.Range("E44").Value = "  -2.30%  "
$ws.Range("D45").Value = "1.795.88"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "'58.01"
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("D47").Value = "0.0₈106"
$ws.Range("E47").Value = "  +4.06%  "
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("D49").Value = "'8.039"
$ws.Range("E49").Value = "  +0.28%  "
$ws.Range("D50").Value = "'0.4346"
$ws.Range("E50").Value = "  -1.28%  "
$ws.Range("D51").Value = "'0.05208"
$ws.Range("E51").Value = "  -1.25%  "
